$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.66920033333333
$ws.Range("H2").Value = 47.00760099999999
$ws.Range("I2").Value = 0.2925937299273087
$ws.Range("J2").Value = 0.2925937299273087
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 55.908252
$ws.Range("N2").Value = 167.724756
$ws.Range("O2").Value = 0.6412441619121594
$ws.Range("P2").Value = 0.6412441619121594
$ws.Range("Q2").Value = 876.0376008744839
$ws.Range("R2").Value = 7884.338407870356
$ws.Range("S2").Value = 0.1876240211279898
$ws.Range("T2").Value = 0.1876240211279898

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.66920033333333
$ws.Range("H3").Value = 47.00760099999999
$ws.Range("I3").Value = 0.2925937299273087
$ws.Range("J3").Value = 0.2925937299273087
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.370676
$ws.Range("N3").Value = 16.112028
$ws.Range("O3").Value = 0.06159939735768789
$ws.Range("P3").Value = 0.06159939735768789
$ws.Range("Q3").Value = 84.15419816942533
$ws.Range("R3").Value = 757.387783524828
$ws.Range("S3").Value = 0.0180235974341603
$ws.Range("T3").Value = 0.0180235974341603

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.66920033333333
$ws.Range("H4").Value = 47.00760099999999
$ws.Range("I4").Value = 0.2925937299273087
$ws.Range("J4").Value = 0.2925937299273087
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 25.90822366666667
$ws.Range("N4").Value = 77.724671
$ws.Range("O4").Value = 0.2971564407301527
$ws.Range("P4").Value = 0.2971564407301527
$ws.Range("Q4").Value = 405.9611469138079
$ws.Range("R4").Value = 3653.65032222427
$ws.Range("S4").Value = 0.0869461113651586
$ws.Range("T4").Value = 0.0869461113651586

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.47676966666667
$ws.Range("H5").Value = 49.43030900000001
$ws.Range("I5").Value = 0.3076736139282969
$ws.Range("J5").Value = 0.3076736139282968
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 55.908252
$ws.Range("N5").Value = 167.724756
$ws.Range("O5").Value = 0.6412441619121594
$ws.Range("P5").Value = 0.6412441619121594
$ws.Range("Q5").Value = 921.1873906699562
$ws.Range("R5").Value = 8290.686516029606
$ws.Range("S5").Value = 0.197293908705936
$ws.Range("T5").Value = 0.197293908705936

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.47676966666667
$ws.Range("H6").Value = 49.43030900000001
$ws.Range("I6").Value = 0.3076736139282969
$ws.Range("J6").Value = 0.3076736139282968
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.370676
$ws.Range("N6").Value = 16.112028
$ws.Range("O6").Value = 0.06159939735768789
$ws.Range("P6").Value = 0.06159939735768789
$ws.Range("Q6").Value = 88.4913914062947
$ws.Range("R6").Value = 796.4225226566523
$ws.Range("S6").Value = 0.01895250920084501
$ws.Range("T6").Value = 0.01895250920084501

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.47676966666667
$ws.Range("H7").Value = 49.43030900000001
$ws.Range("I7").Value = 0.3076736139282969
$ws.Range("J7").Value = 0.3076736139282968
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 25.90822366666667
$ws.Range("N7").Value = 77.724671
$ws.Range("O7").Value = 0.2971564407301527
$ws.Range("P7").Value = 0.2971564407301527
$ws.Range("Q7").Value = 426.8838338281489
$ws.Range("R7").Value = 3841.95450445334
$ws.Range("S7").Value = 0.09142719602151583
$ws.Range("T7").Value = 0.09142719602151582

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 21.406785
$ws.Range("H8").Value = 64.220355
$ws.Range("I8").Value = 0.3997326561443945
$ws.Range("J8").Value = 0.3997326561443944
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 55.908252
$ws.Range("N8").Value = 167.724756
$ws.Range("O8").Value = 0.6412441619121594
$ws.Range("P8").Value = 0.6412441619121594
$ws.Range("Q8").Value = 1196.81593028982
$ws.Range("R8").Value = 10771.34337260838
$ws.Range("S8").Value = 0.2563262320782336
$ws.Range("T8").Value = 0.2563262320782336

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 21.406785
$ws.Range("H9").Value = 64.220355
$ws.Range("I9").Value = 0.3997326561443945
$ws.Range("J9").Value = 0.3997326561443944
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.370676
$ws.Range("N9").Value = 16.112028
$ws.Range("O9").Value = 0.06159939735768789
$ws.Range("P9").Value = 0.06159939735768789
$ws.Range("Q9").Value = 114.96890643666
$ws.Range("R9").Value = 1034.72015792994
$ws.Range("S9").Value = 0.02462329072268257
$ws.Range("T9").Value = 0.02462329072268257

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 21.406785
$ws.Range("H10").Value = 64.220355
$ws.Range("I10").Value = 0.3997326561443945
$ws.Range("J10").Value = 0.3997326561443944
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 25.90822366666667
$ws.Range("N10").Value = 77.724671
$ws.Range("O10").Value = 0.2971564407301527
$ws.Range("P10").Value = 0.2971564407301527
$ws.Range("Q10").Value = 554.611773764245
$ws.Range("R10").Value = 4991.505963878205
$ws.Range("S10").Value = 0.1187831333434783
$ws.Range("T10").Value = 0.1187831333434782
